$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")

# Rename the "capacity_market_active" setting to "capacity_remuneration_mechanism"
# and turn its value column from a boolean switch into a chosen mechanism string,
# adding an explanatory comment of the allowed options (CRM support).
$ws.Range("A44").Value = "capacity_remuneration_mechanism"
$ws.Range("C44").Value = '"capacity_market", "strategic_reserve", "strategic_reserve_ger", "strategic_reserve_swe", "forward_capacity_market"'
$ws.Range("B44").Value = "strategic_reserve_ger"

# Narrow column C now that its content is a shorter list instead of the old wide text
$ws.Columns.Item(3).ColumnWidth = 62.8

# Update the saved selection/scroll position on the sheet
$null = $ws.Range("H34").Select()
